$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, E, F, G, H, K, L, M across rows 2-25
$colB = @(0.72996587998216, 0.7015787694087408, 0.6845469246810012, 0.6777067096255962, 0.6765769708721336, 0.6844542681918711, 0.7200955949786021, 0.7931358398373902, 0.8487110012825099, 0.8744079914804104, 0.8841983169077139, 0.8820871570034114, 0.8752122581724677, 0.8710089144556434, 0.8470399769593655, 0.8324420148414049, 0.8240847878159343, 0.8212619057460984, 0.8339919449079503, 0.8772299721912589, 0.9058347948135292, 0.8905362852758003, 0.8332911114337094, 0.7730402146849258)
$colC = @(0.1260741705134976, 0.1249385211811003, 0.1242236898776916, 0.1239279735702823, 0.1238786031521641, 0.1242197196352883, 0.1256862343411598, 0.1284233112840738, 0.1303503566079627, 0.1312089661045093, 0.1315315187619461, 0.1314621661063171, 0.1312355544598702, 0.1310964119389411, 0.1302938825099318, 0.1297969474287797, 0.1295094276618087, 0.1294117869589044, 0.1298500224877515, 0.1313021858768266, 0.1322361992991929, 0.1317390746148845, 0.1298260329405352, 0.1276976419417863)
$colE = @(0.1131151126661525, 0.1132122997500495, 0.1133083146848186, 0.1133565946677155, 0.1133651647815075, 0.1133089287252762, 0.1131410887035518, 0.1130997393740181, 0.1132441465048721, 0.1133476641460653, 0.1133922902675373, 0.1133824381087152, 0.1133512268889731, 0.1133328153839024, 0.1132381414196004, 0.1131897408359563, 0.1131654632540258, 0.113157855313867, 0.1131945246900532, 0.1133602472184414, 0.113500178691492, 0.1134226055123015, 0.1131923508555026, 0.1130801840602693)
$colF = @(0.4443680307746263, 0.387822817061874, 0.3531389305168915, 0.3390132514313251, 0.336668177824194, 0.3529483938344953, 0.4248636149813336, 0.5661985755041457, 0.6702781546542269, 0.7176906081379002, 0.7356546913071611, 0.7317853510981394, 0.7191683204515869, 0.7114413442032514, 0.6671810134426437, 0.6400460337125793, 0.6244449056556647, 0.619163680173358, 0.642933953830422, 0.7228739723491628, 0.7751780083420101, 0.7472568307830727, 0.6416283278902171, 0.5279251897347166)
$colG = @(1.881139910082354, 1.878022025669864, 1.876923551168503, 1.876680803271583, 1.876652861791257, 1.876919448173552, 1.879895346486094, 1.892219809543747, 1.905253866663713, 1.912052720256895, 1.914752680743902, 1.914165614203569, 1.912272332867332, 1.911128981783691, 1.904827072766693, 1.901184018898078, 1.899170461107815, 1.898502749269852, 1.901563356621153, 1.912825029983708, 1.920916163018347, 1.916530777840336, 1.901391606110224, 1.888188613390327)
$colH = @(1.62608349991342, 1.629574383295832, 1.632271618605628, 1.633510005374532, 1.633724047806325, 1.632287756172857, 1.627172197010879, 1.621536778939415, 1.620080420160264, 1.62000164375101, 1.620055809977856, 1.620040407591461, 1.620004416279556, 1.619993311187443, 1.620097319416487, 1.620310672907408, 1.620488327987374, 1.620557913163566, 1.620282274578727, 1.620012707679649, 1.620326156836541, 1.620114043085778, 1.620294942149769, 1.622590260580054)
$colK = @(0.3524533962567205, 0.3260699766826605, 0.3100425494185259, 0.30355471395319, 0.3024800452238736, 0.3099548759250155, 0.3433207814575212, 0.4101119950857992, 0.4600123503855116, 0.4828935470738145, 0.4915840288818174, 0.4897112307085365, 0.4836080009069121, 0.4798729595428313, 0.4585206410936848, 0.4454679723367008, 0.4379774965668162, 0.4354442887448329, 0.4468556855506733, 0.4853999669223867, 0.5107415789486538, 0.4972025595424157, 0.4462282577898407, 0.3918976081790504)
$colL = @(0.1949625722139459, 0.1925609050790555, 0.1911798392513191, 0.1906406110580576, 0.1905524971959309, 0.1911724715556886, 0.1941150795467692, 0.2006267630097298, 0.2058619868255818, 0.2083414654051836, 0.209294439175693, 0.2090885749566809, 0.2084195857713667, 0.2080116391190074, 0.2057019149774248, 0.2043100354701437, 0.2035186847430026, 0.2032523318938075, 0.2044572492853547, 0.2086157031090039, 0.2114153652275519, 0.2099136533354908, 0.2043906663215722, 0.198785908276264)
$colM = @(0.194132664201959, 0.1887435681158003, 0.1855385113821413, 0.184258611349815, 0.1840476685420214, 0.1855211440837188, 0.192252983813372, 0.2062762401527003, 0.2170790081025444, 0.2221018332932729, 0.2240194181312347, 0.2236057413724737, 0.2222592828800387, 0.2214365613079039, 0.2167529342001586, 0.2139074468703157, 0.2122810228655894, 0.2117321018105116, 0.2142092960767101, 0.2226543488746415, 0.2282642862156621, 0.2252618879486903, 0.2140728005025707, 0.202394701055038)

$colsMap = @{
    2 = $colB
    3 = $colC
    5 = $colE
    6 = $colF
    7 = $colG
    8 = $colH
    11 = $colK
    12 = $colL
    13 = $colM
}

$startRow = 2
foreach ($colIndex in $colsMap.Keys) {
    $vals = $colsMap[$colIndex]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($startRow + $i, $colIndex).Value = $vals[$i]
    }
}

Write-Host "Applied 380 kV case updates"
